$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2157.5264
$ws.Range("I6").Value = 723.6
$ws.Range("K6").Value = 2170.8
$ws.Range("M6").Value = -2058.8
$ws.Range("H18").Value = 2119.6
$ws.Range("I18").Value = 2119.6
$ws.Range("K18").Value = 2119.6
$ws.Range("M18").Value = -1835.6
$ws.Range("H80").Value = 1033814.75
$ws.Range("I80").Value = 2066778
$ws.Range("J80").Value = 851.5454999999999
$ws.Range("K80").Value = 6200334
$ws.Range("L80").Value = 2554.6365
$ws.Range("M80").Value = -6199336
$ws.Range("N80").Value = -4550.6365
$ws.Range("H83").Value = 1033814.75
$ws.Range("I83").Value = 2066778
$ws.Range("J83").Value = 851.5454999999999
$ws.Range("K83").Value = 18601002
$ws.Range("L83").Value = 7663.9095
$ws.Range("M83").Value = -18596010
$ws.Range("N83").Value = -17647.9095
$ws.Range("H100").Value = 3787.7856
$ws.Range("I100").Value = 4066.0908
$ws.Range("K100").Value = 4066.0908
$ws.Range("M100").Value = -3525.0908
$ws.Range("H137").Value = 2412.07
$ws.Range("I137").Value = 1831
$ws.Range("J137").Value = 2476.6333
$ws.Range("K137").Value = 5493
$ws.Range("L137").Value = 7429.8999
$ws.Range("M137").Value = -2943
$ws.Range("N137").Value = -12529.8999
$ws.Range("H141").Value = 3120.5
$ws.Range("I141").Value = 3145.2666
$ws.Range("J141").Value = 2996.6667
$ws.Range("K141").Value = 9435.799800000001
$ws.Range("L141").Value = 8990.000100000001
$ws.Range("M141").Value = -4255.799800000001
$ws.Range("N141").Value = -19350.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1178495.8
$ws.Range("I61").Value = 668761.9
$ws.Range("K61").Value = 668761.9
$ws.Range("M61").Value = -668549.9
$ws.Range("H74").Value = 2744.5833
$ws.Range("I74").Value = 2908.6843
$ws.Range("J74").Value = 2121
$ws.Range("K74").Value = 2908.6843
$ws.Range("L74").Value = 2121
$ws.Range("M74").Value = -2034.6843
$ws.Range("N74").Value = -3869
$ws.Range("H77").Value = 2744.5833
$ws.Range("I77").Value = 2908.6843
$ws.Range("J77").Value = 2121
$ws.Range("K77").Value = 14543.4215
$ws.Range("L77").Value = 10605
$ws.Range("M77").Value = -10175.4215
$ws.Range("N77").Value = -19341
$ws.Range("H132").Value = 742486.9
$ws.Range("I132").Value = 477506.94
$ws.Range("K132").Value = 1432520.82
$ws.Range("M132").Value = -1429990.82
$ws.Range("H136").Value = 1178495.8
$ws.Range("I136").Value = 668761.9
$ws.Range("K136").Value = 2006285.7
$ws.Range("M136").Value = -2003735.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 74000
$ws.Range("J35").Value = 74000
$ws.Range("L35").Value = 74000
$ws.Range("N35").Value = -74620
$ws.Range("H86").Value = 1406.238
$ws.Range("I86").Value = 1425.4706
$ws.Range("K86").Value = 1425.4706
$ws.Range("M86").Value = -302.4706000000001
$ws.Range("H89").Value = 1406.238
$ws.Range("I89").Value = 1425.4706
$ws.Range("K89").Value = 7127.353000000001
$ws.Range("M89").Value = -1511.353000000001
$ws.Range("H134").Value = 3751.2222
$ws.Range("J134").Value = 4159.5
$ws.Range("L134").Value = 12478.5
$ws.Range("N134").Value = -17548.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1879.1923
$ws.Range("I58").Value = 1280.0834
$ws.Range("J58").Value = 2392.7144
$ws.Range("K58").Value = 1280.0834
$ws.Range("L58").Value = 2392.7144
$ws.Range("M58").Value = -1077.0834
$ws.Range("N58").Value = -2798.7144
$ws.Range("H62").Value = 8567.083000000001
$ws.Range("I62").Value = 8436.817999999999
$ws.Range("K62").Value = 8436.817999999999
$ws.Range("M62").Value = -7812.817999999999
$ws.Range("H65").Value = 8567.083000000001
$ws.Range("I65").Value = 8436.817999999999
$ws.Range("K65").Value = 42184.09
$ws.Range("M65").Value = -39064.09
$ws.Range("H122").Value = 2730.9375
$ws.Range("I122").Value = 1707.3846
$ws.Range("J122").Value = 7166.3335
$ws.Range("K122").Value = 5122.1538
$ws.Range("L122").Value = 21499.0005
$ws.Range("M122").Value = -2672.1538
$ws.Range("N122").Value = -26399.0005
$ws.Range("H136").Value = 1879.1923
$ws.Range("I136").Value = 1280.0834
$ws.Range("J136").Value = 2392.7144
$ws.Range("K136").Value = 3840.2502
$ws.Range("L136").Value = 7178.1432
$ws.Range("M136").Value = -1290.2502
$ws.Range("N136").Value = -12278.1432
$ws.Range("H141").Value = 71999.75
$ws.Range("J141").Value = 84833
$ws.Range("L141").Value = 84833
$ws.Range("N141").Value = -95193

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 391.83334
$ws.Range("I10").Value = 330.25
$ws.Range("K10").Value = 990.75
$ws.Range("M10").Value = -851.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6163.3335
$ws.Range("I12").Value = 1745
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 1745
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -1605
$ws.Range("N12").Value = -15280
$ws.Range("H111").Value = 41999
$ws.Range("J111").Value = 41999
$ws.Range("L111").Value = 41999
$ws.Range("N111").Value = -48133

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25003166
$ws.Range("I7").Value = 38464176
$ws.Range("K7").Value = 38464176
$ws.Range("M7").Value = -38464064
$ws.Range("H40").Value = 3166.6875
$ws.Range("I40").Value = 2606.5454
$ws.Range("K40").Value = 2606.5454
$ws.Range("M40").Value = -2470.5454
$ws.Range("H46").Value = 2962.0625
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 3199.5
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 3199.5
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -3575.5
$ws.Range("H82").Value = 3490.182
$ws.Range("I82").Value = 1248
$ws.Range("J82").Value = 4771.4287
$ws.Range("K82").Value = 1248
$ws.Range("L82").Value = 4771.4287
$ws.Range("M82").Value = -887
$ws.Range("N82").Value = -5493.4287
$ws.Range("H85").Value = 3490.182
$ws.Range("I85").Value = 1248
$ws.Range("J85").Value = 4771.4287
$ws.Range("K85").Value = 1248
$ws.Range("L85").Value = 4771.4287
$ws.Range("M85").Value = 0
$ws.Range("N85").Value = -7267.4287
$ws.Range("H126").Value = 25003166
$ws.Range("I126").Value = 38464176
$ws.Range("K126").Value = 115392528
$ws.Range("M126").Value = -115390058
$ws.Range("H132").Value = 5727.143
$ws.Range("I132").Value = 3306.2354
$ws.Range("J132").Value = 9468.546
$ws.Range("K132").Value = 9918.706200000001
$ws.Range("L132").Value = 28405.638
$ws.Range("M132").Value = -7388.706200000001
$ws.Range("N132").Value = -33465.638
$ws.Range("H136").Value = 7463.1816
$ws.Range("I136").Value = 4057.5715
$ws.Range("K136").Value = 12172.7145
$ws.Range("M136").Value = -9622.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1607.579
$ws.Range("I126").Value = 1551.1765
$ws.Range("K126").Value = 4653.529500000001
$ws.Range("M126").Value = -2183.529500000001
$ws.Range("H136").Value = 4365.2104
$ws.Range("I136").Value = 3540.742
$ws.Range("J136").Value = 8016.4287
$ws.Range("K136").Value = 10622.226
$ws.Range("L136").Value = 24049.2861
$ws.Range("M136").Value = -8072.226000000001
$ws.Range("N136").Value = -29149.2861
